# Auto-generated Excel COM-interop script
# Applies targeted cell value updates/additions/deletions across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets as described by the source diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 827
$ws.Range("I19").Value = 715.3333
$ws.Range("K19").Value = 715.3333
$ws.Range("M19").Value = -540.3333
$ws.Range("H40").Value = 5917.1665
$ws.Range("I40").Value = 6833
$ws.Range("J40").Value = 5001.3335
$ws.Range("K40").Value = 6833
$ws.Range("L40").Value = 5001.3335
$ws.Range("M40").Value = -6658
$ws.Range("N40").Value = -5351.3335
$ws.Range("H86").Value = 8941.799999999999
$ws.Range("I86").Value = 8659.272000000001
$ws.Range("K86").Value = 8659.272000000001
$ws.Range("M86").Value = -7536.272000000001
$ws.Range("H89").Value = 8941.799999999999
$ws.Range("I89").Value = 8659.272000000001
$ws.Range("K89").Value = 43296.36
$ws.Range("M89").Value = -37680.36
$ws.Range("H111").Value = 4477.706
$ws.Range("J111").Value = 1035.25
$ws.Range("L111").Value = 3105.75
$ws.Range("N111").Value = -9239.75
$ws.Range("H125").Value = 945.73334
$ws.Range("I125").Value = 839.2
$ws.Range("J125").Value = 999
$ws.Range("K125").Value = 7552.8
$ws.Range("L125").Value = 8991
$ws.Range("M125").Value = -5092.8
$ws.Range("N125").Value = -13911
$ws.Range("H137").Value = 2001898.9
$ws.Range("I137").Value = 2633537.2
$ws.Range("J137").Value = 1710.5
$ws.Range("K137").Value = 7900611.600000001
$ws.Range("L137").Value = 5131.5
$ws.Range("M137").Value = -7898061.600000001
$ws.Range("N137").Value = -10231.5
$ws.Range("H138").Value = 739.5
$ws.Range("I138").Value = 739.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2218.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2921.5
$ws.Range("N138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 49039.11
$ws.Range("I45").Value = 49039.11
$ws.Range("K45").Value = 49039.11
$ws.Range("M45").Value = -48662.11
$ws.Range("H61").Value = 1741.8379
$ws.Range("I61").Value = 1200.75
$ws.Range("K61").Value = 1200.75
$ws.Range("M61").Value = -988.75
$ws.Range("H132").Value = 2335.3333
$ws.Range("I132").Value = 1503
$ws.Range("K132").Value = 4509
$ws.Range("M132").Value = -1979
$ws.Range("H136").Value = 1741.8379
$ws.Range("I136").Value = 1200.75
$ws.Range("K136").Value = 3602.25
$ws.Range("M136").Value = -1052.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3334.75
$ws.Range("I99").Value = 2046.25
$ws.Range("J99").Value = 4623.25
$ws.Range("K99").Value = 2046.25
$ws.Range("L99").Value = 4623.25
$ws.Range("M99").Value = -548.25
$ws.Range("N99").Value = -7619.25
$ws.Range("H134").Value = 2969.7693
$ws.Range("I134").Value = 2338.5
$ws.Range("J134").Value = 3979.8
$ws.Range("K134").Value = 7015.5
$ws.Range("L134").Value = 11939.4
$ws.Range("M134").Value = -4480.5
$ws.Range("N134").Value = -17009.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6948420
$ws.Range("J31").Value = 17861518
$ws.Range("L31").Value = 17861518
$ws.Range("N31").Value = -17862108
$ws.Range("H34").Value = 6948420
$ws.Range("J34").Value = 17861518
$ws.Range("L34").Value = 17861518
$ws.Range("N34").Value = -17861922
$ws.Range("H62").Value = 397
$ws.Range("I62").Value = 397
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 397
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 227
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 397
$ws.Range("I65").Value = 397
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1985
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 1135
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 4604.778
$ws.Range("I99").Value = 3906
$ws.Range("J99").Value = 5163.8
$ws.Range("K99").Value = 3906
$ws.Range("L99").Value = 5163.8
$ws.Range("M99").Value = -2408
$ws.Range("N99").Value = -8159.8
$ws.Range("H126").Value = 4604.778
$ws.Range("I126").Value = 3906
$ws.Range("J126").Value = 5163.8
$ws.Range("K126").Value = 11718
$ws.Range("L126").Value = 15491.4
$ws.Range("M126").Value = -9248
$ws.Range("N126").Value = -20431.4
$ws.Range("H134").Value = 4850.76
$ws.Range("I134").Value = 5065.2856
$ws.Range("J134").Value = 3724.5
$ws.Range("K134").Value = 15195.8568
$ws.Range("L134").Value = 11173.5
$ws.Range("M134").Value = -12660.8568
$ws.Range("N134").Value = -16243.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 6800
$ws.Range("I99").Value = 6800
$ws.Range("K99").Value = 20400
$ws.Range("M99").Value = -18154
$ws.Range("H131").Value = 1790234.8
$ws.Range("I131").Value = 49507.5
$ws.Range("K131").Value = 148522.5
$ws.Range("M131").Value = -143482.5
$ws.Range("H132").Value = 1709.8572
$ws.Range("I132").Value = 1274.25
$ws.Range("K132").Value = 11468.25
$ws.Range("M132").Value = -8938.25
$ws.Range("H133").Value = 3666.5833
$ws.Range("H134").Value = 2490.6
$ws.Range("J134").Value = 5000
$ws.Range("L134").Value = 15000
$ws.Range("N134").Value = -25140
$ws.Range("H137").Value = 2266.6667
$ws.Range("I137").Value = 2400
$ws.Range("K137").Value = 7200
$ws.Range("M137").Value = -2100

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4999.1377
$ws.Range("J102").Value = 5459.72
$ws.Range("L102").Value = 5459.72
$ws.Range("N102").Value = -8703.720000000001
$ws.Range("H132").Value = 2069.3333
$ws.Range("I132").Value = 2004.8
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 6014.4
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -3484.4
$ws.Range("N132").Value = -11510

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 48500
$ws.Range("I40").Value = 48500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 48500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -48364
$ws.Range("N40").ClearContents()
$ws.Range("H48").Value = 31666.666
$ws.Range("J48").Value = 35000
$ws.Range("L48").Value = 35000
$ws.Range("N48").Value = -36322
$ws.Range("H68").Value = 1179.8
$ws.Range("I68").Value = 1133
$ws.Range("J68").Value = 1250
$ws.Range("K68").Value = 1133
$ws.Range("L68").Value = 1250
$ws.Range("M68").Value = -384
$ws.Range("N68").Value = -2748
$ws.Range("H71").Value = 1179.8
$ws.Range("I71").Value = 1133
$ws.Range("J71").Value = 1250
$ws.Range("K71").Value = 5665
$ws.Range("L71").Value = 6250
$ws.Range("M71").Value = -1921
$ws.Range("N71").Value = -13738
$ws.Range("H93").Value = 2820.4443
$ws.Range("I93").Value = 2675.5557
$ws.Range("K93").Value = 2675.5557
$ws.Range("M93").Value = -1427.5557
$ws.Range("H122").Value = 6312.7144
$ws.Range("I122").Value = 2651.5
$ws.Range("J122").Value = 7777.2
$ws.Range("K122").Value = 7954.5
$ws.Range("L122").Value = 23331.6
$ws.Range("M122").Value = -5504.5
$ws.Range("N122").Value = -28231.6
$ws.Range("H132").Value = 4551.5557
$ws.Range("I132").Value = 4712.727
$ws.Range("K132").Value = 14138.181
$ws.Range("M132").Value = -11608.181
$ws.Range("H133").Value = 129998.5
$ws.Range("J133").Value = 129998.5
$ws.Range("L133").Value = 129998.5
$ws.Range("N133").Value = -135058.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22729436
$ws.Range("I122").Value = 2365.111
$ws.Range("K122").Value = 7095.333
$ws.Range("M122").Value = -4645.333
$ws.Range("H126").Value = 7825.1904
$ws.Range("I126").Value = 8122.6313
$ws.Range("K126").Value = 24367.8939
$ws.Range("M126").Value = -21897.8939
$ws.Range("H132").Value = 3537.25
$ws.Range("I132").Value = 3610.4
$ws.Range("K132").Value = 10831.2
$ws.Range("M132").Value = -8301.200000000001

